# Adds two new EIR model runs (Alt1_v2 and Alt2_v1) to the all_runs sheet,
# each contributing a 2035 and a 2050 row - matching the commit message
# "added EIR alt1_v2 and alt2_v1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# Find the last used row in column A so we append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # -4162 = xlUp

# New rows to append, in order.
$newRows = @(
    @{ Project="RTP2021"; Year=2035; Directory="2035_TM152_EIR_Alt1_02"; RunSet="EIR"; Category="Alt1"; UrbansimPath='"EIR runs\Alt1 (s26) runs\Alt1_v2"'; RunId="run373"; Status="current" },
    @{ Project="RTP2021"; Year=2050; Directory="2050_TM152_EIR_Alt1_02"; RunSet="EIR"; Category="Alt1"; UrbansimPath='"EIR runs\Alt1 (s26) runs\Alt1_v2"'; RunId="run373"; Status="current" },
    @{ Project="RTP2021"; Year=2035; Directory="2035_TM152_EIR_Alt2_01"; RunSet="EIR"; Category="Alt2"; UrbansimPath='"EIR runs\Alt2 (s28) runs\Alt2_v1"'; RunId="run374"; Status="current" },
    @{ Project="RTP2021"; Year=2050; Directory="2050_TM152_EIR_Alt2_01"; RunSet="EIR"; Category="Alt2"; UrbansimPath='"EIR runs\Alt2 (s28) runs\Alt2_v1"'; RunId="run374"; Status="current" }
)

$templateRow = $lastRow

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $lastRow + 1 + $i

    # Copy formatting (styles) from the last existing row down onto the new row.
    $ws.Range("A$templateRow`:H$templateRow").Copy() | Out-Null
    $ws.Range("A$targetRow`:H$targetRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = 0

    $row = $newRows[$i]

    $ws.Cells.Item($targetRow, 1).Value = $row.Project
    $ws.Cells.Item($targetRow, 2).Value = $row.Year
    $ws.Cells.Item($targetRow, 3).Value = $row.Directory
    $ws.Cells.Item($targetRow, 4).Value = $row.RunSet
    $ws.Cells.Item($targetRow, 5).Value = $row.Category
    $ws.Cells.Item($targetRow, 6).Value = $row.UrbansimPath
    $ws.Cells.Item($targetRow, 7).Value = $row.RunId
    $ws.Cells.Item($targetRow, 8).Value = $row.Status
}

# Leave the sheet selection on the newly-added rows, as it was when the
# author saved the workbook after appending them.
$ws.Activate()
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $newRows.Count
$ws.Rows("$firstNewRow`:$lastNewRow").Select() | Out-Null
